$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44643
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 28000
$ws.Range("P2").Value = 29000
$ws.Range("S2").Value = 1450

# Row 3
$ws.Range("D3").Value = 44650
$ws.Range("N3").Value = 31000
$ws.Range("O3").Value = 32000
$ws.Range("P3").Value = 31500
$ws.Range("S3").Value = 1575

# Row 4
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 29000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 29500
$ws.Range("S4").Value = 1475

# Row 5
$ws.Range("D5").Value = 44671
$ws.Range("M5").Value = 200

# Row 6
$ws.Range("D6").Value = 44636
$ws.Range("L6").Value = "Primera"
